# Applies the "456a3b4" data refresh to 江西-漫展信息.xlsx:
#   - bumps the "想去人数" (interest count) figures in column F on both the
#     "展览" sheet and the combined "全部类型" sheet
#   - refreshes one cover-image URL (row 22 / 23 respectively)
#   - inserts a newly-scraped con ("南昌·LY-COSPLAY大会...") into the
#     "展览" sheet as row 35 (pushing the old rows 35-37 down to 36-38,
#     and renumbering their index column) - it was already present as
#     row 36 of "全部类型", where only its F (interest count) value changes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibition list)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 5557
$ws1.Range("F5").Value = 60
$ws1.Range("F7").Value = 662
$ws1.Range("F8").Value = 644
$ws1.Range("F9").Value = 9
$ws1.Range("F12").Value = 1553
$ws1.Range("F13").Value = 5109
$ws1.Range("F15").Value = 247
$ws1.Range("F16").Value = 216
$ws1.Range("F17").Value = 33
$ws1.Range("F19").Value = 112
$ws1.Range("F20").Value = 4412
$ws1.Range("F21").Value = 212
$ws1.Range("F22").Value = 1161
$ws1.Range("F23").Value = 121
$ws1.Range("F24").Value = 68
$ws1.Range("F25").Value = 212
$ws1.Range("F26").Value = 58
$ws1.Range("F27").Value = 176
$ws1.Range("F29").Value = 149
$ws1.Range("F31").Value = 346

# Refreshed cover image for the "New World" entry (row 22)
$ws1.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202404/eECmAU3Q1713511863782.jpeg"

# Insert a fresh row at 35 for the new con, pushing the former rows
# 35-37 down to 36-38 (dimension grows from A1:I37 to A1:I38).
$ws1.Rows.Item(35).Insert()

# Copy column A's formatting (bold + border + centered) onto the new
# row's index cell, then give it its value.
$ws1.Range("A2").Copy()
$ws1.Range("A35").PasteSpecial(-4122)
$ws1.Range("A35").Value = 34

# Column B holds literal text dates ("2024-06-10"), not real dates, so
# force text formatting before assigning, then drop the leftover
# number-format style so the cell matches its siblings again.
$ws1.Range("B35").NumberFormat = "@"
$ws1.Range("B35").Value = "2024-06-10"
$ws1.Range("B35").ClearFormats()

$ws1.Range("C35").Value = "南昌·LY-COSPLAY大会X运动番PRO2.0（非ONLY）"
$ws1.Range("D35").Value = "青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK"
$ws1.Range("E35").Value = "2024.06.10 10:00-06.10 17:00"
$ws1.Range("F35").Value = 8
$ws1.Range("G35").Value = 30
$ws1.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=84575"
$ws1.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202404/ScwkijwU1713428452963.jpeg"

# The index column (A) is a plain literal sequence (0-based row number),
# not a formula, so the rows that got pushed down need to be renumbered
# by hand to keep counting up from the new row.
$ws1.Range("A36").Value = 35
$ws1.Range("A37").Value = 36
$ws1.Range("A38").Value = 37

# ---------------------------------------------------------------------
# Sheet "全部类型" (combined list - already contains the new con as row
# 36, only its numbers/links need the same refresh as above)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 5557
$ws4.Range("F6").Value = 60
$ws4.Range("F8").Value = 662
$ws4.Range("F9").Value = 644
$ws4.Range("F10").Value = 9
$ws4.Range("F13").Value = 1553
$ws4.Range("F14").Value = 5109
$ws4.Range("F16").Value = 247
$ws4.Range("F17").Value = 216
$ws4.Range("F18").Value = 33
$ws4.Range("F20").Value = 112
$ws4.Range("F21").Value = 4412
$ws4.Range("F22").Value = 212
$ws4.Range("F23").Value = 1161
$ws4.Range("F24").Value = 121
$ws4.Range("F25").Value = 68
$ws4.Range("F26").Value = 212
$ws4.Range("F27").Value = 58
$ws4.Range("F28").Value = 176
$ws4.Range("F30").Value = 149
$ws4.Range("F32").Value = 346
$ws4.Range("F36").Value = 8

$ws4.Range("I23").Value = "//i1.hdslb.com/bfs/openplatform/202404/eECmAU3Q1713511863782.jpeg"
